$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (prevents "44.520.09"-style strings from being parsed as numbers/dates),
    # then reset the style to Normal so no stray quotePrefix/number-format
    # style gets attached to the cell.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "44.457.28"
Set-TextValue $ws.Range("E2") "  +3.84%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.284.75"
Set-TextValue $ws.Range("E3") "  +2.88%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "320.67"
Set-TextValue $ws.Range("E5") "  +1.38%  "

# Row 6
Set-TextValue $ws.Range("D6") "107.71"
Set-TextValue $ws.Range("E6") "  +7.81%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  -0.02%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.577"
Set-TextValue $ws.Range("E9") "  +2.53%  "

# Row 10
Set-TextValue $ws.Range("D10") "39.06"
Set-TextValue $ws.Range("E10") "  +5.04%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +2.21%  "

# Row 12
Set-TextValue $ws.Range("D12") "7.98"
Set-TextValue $ws.Range("E12") "  +2.01%  "

# Row 13
Set-TextValue $ws.Range("E13") "  +1.61%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.890"

# Row 15
Set-TextValue $ws.Range("D15") "2.633.89"
Set-TextValue $ws.Range("E15") "  +3.01%  "

# Row 16
Set-TextValue $ws.Range("D16") "14.73"
Set-TextValue $ws.Range("E16") "  +3.33%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.340.05"
Set-TextValue $ws.Range("E17") "  +5.49%  "

# Row 18
Set-TextValue $ws.Range("D18") "44.360.17"
Set-TextValue $ws.Range("E18") "  +3.71%  "

# Row 19
Set-TextValue $ws.Range("D19") "14.22"
Set-TextValue $ws.Range("E19") "  -8.16%  "

# Row 20
Set-TextValue $ws.Range("E20") "  +4.33%  "

# Row 21
Set-TextValue $ws.Range("E21") "  +2.24%  "

# Row 22
Set-TextValue $ws.Range("D22") "66.77"
Set-TextValue $ws.Range("E22") "  +2.18%  "

# Row 23
Set-TextValue $ws.Range("D23") "3.23"
Set-TextValue $ws.Range("E23") "  +2.43%  "

# Row 24
Set-TextValue $ws.Range("D24") "239.93"
Set-TextValue $ws.Range("E24") "  +1.44%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.22"
Set-TextValue $ws.Range("E25") "  +4.05%  "

# Row 26
Set-TextValue $ws.Range("E26") "  -0.30%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.31"
Set-TextValue $ws.Range("E27") "  +1.87%  "

# Row 28
Set-TextValue $ws.Range("D28") "39.08"
Set-TextValue $ws.Range("E28") "  +13.96%  "

# Row 29
Set-TextValue $ws.Range("E29") "  +0.67%  "

# Row 30
Set-TextValue $ws.Range("D30") "6.58"
Set-TextValue $ws.Range("E30") "  +3.64%  "

# Row 31
Set-TextValue $ws.Range("D31") "20.82"
Set-TextValue $ws.Range("E31") "  +1.45%  "

# Row 32
Set-TextValue $ws.Range("D32") "163.44"
Set-TextValue $ws.Range("E32") "  +4.32%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0890"
Set-TextValue $ws.Range("E33") "  +0.86%  "

# Row 34
Set-TextValue $ws.Range("D34") "2.74"
Set-TextValue $ws.Range("E34") "  -1.19%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "2.08"
Set-TextValue $ws.Range("E35") "  +5.40%  "

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D36") "3.28"
Set-TextValue $ws.Range("E36") "  +2.40%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.117"
Set-TextValue $ws.Range("E37") "  +13.40%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D38") "0.122"
Set-TextValue $ws.Range("E38") "  -0.74%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D39") "4.00"
Set-TextValue $ws.Range("E39") "  +4.69%  "

# Row 40
Set-TextValue $ws.Range("D40") "4.50"
Set-TextValue $ws.Range("E40") "  +1.90%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +1.75%  "

# Row 42
Set-TextValue $ws.Range("D42") "15.48"
Set-TextValue $ws.Range("E42") "  +26.21%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.00"
Set-TextValue $ws.Range("E43") "  +0.20%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.782.67"
Set-TextValue $ws.Range("E44") "  -7.07%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.209"
Set-TextValue $ws.Range("E45") "  +0.96%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D46") "86.73"
Set-TextValue $ws.Range("E46") "  -2.83%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D47") "5.48"
Set-TextValue $ws.Range("E47") "  +1.60%  "

# Row 48
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D48") "60.72"
Set-TextValue $ws.Range("E48") "  +0.10%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D49") "75.76"
Set-TextValue $ws.Range("E49") "  +0.97%  "

# Row 50
Set-TextValue $ws.Range("D50") "8.80"
Set-TextValue $ws.Range("E50") "  +2.05%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.71"
Set-TextValue $ws.Range("E51") "  +6.32%  "
